$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.722.40'
$ws.Range("E2").Value = '  +3.05%  '
$ws.Range("D3").Value = '1.790.33'
$ws.Range("E3").Value = '  +0.82%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '223.18'
$ws.Range("E5").Value = '  -0.49%  '
$ws.Range("E6").Value = '  -0.19%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '32.53'
$ws.Range("E8").Value = '  +7.90%  '
$ws.Range("E9").Value = '  +1.23%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0689'
$ws.Range("E10").Value = '  +4.30%  '
$ws.Range("E11").Value = '  +1.52%  '
$ws.Range("D12").Value = '2.047.41'
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.813.28'
$ws.Range("E13").Value = '  +2.11%  '
$ws.Range("B14").Value = 'Chainlink'
$ws.Range("C14").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '10.99'
$ws.Range("E14").Value = '  +9.59%  '
$ws.Range("D15").Value = '34.721.77'
$ws.Range("E15").Value = '  +3.11%  '
$ws.Range("E16").Value = '  +1.32%  '
$ws.Range("E18").Value = '  +0.31%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '253.40'
$ws.Range("E19").Value = '  +1.49%  '
$ws.Range("D20").Value = '0.0₃0788'
$ws.Range("E20").Value = '  +7.07%  '
$ws.Range("E21").Value = '  -0.08%  '
$ws.Range("E22").Value = '  +2.07%  '
$ws.Range("E23").Value = '  +0.58%  '
$ws.Range("E24").Value = '  +0.13%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '158.27'
$ws.Range("E25").Value = '  -0.11%  '
$ws.Range("E26").Value = '  -0.07%  '
$ws.Range("E27").Value = '  +1.59%  '
$ws.Range("E28").Value = '  +0.16%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("E29").Value = '  -0.03%  '
$ws.Range("E30").Value = '  +0.46%  '
$ws.Range("E31").Value = '  -1.19%  '
$ws.Range("E32").Value = '  +0.13%  '
$ws.Range("E33").Value = '  +0.60%  '
$ws.Range("E34").Value = '  +2.35%  '
$ws.Range("D35").Value = '1.430.99'
$ws.Range("E35").Value = '  -3.39%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.06'
$ws.Range("E36").Value = '  -0.77%  '
$ws.Range("B37").Value = 'ImmutableX'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.632'
$ws.Range("E37").Value = '  +0.90%  '
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0189'
$ws.Range("E38").Value = '  +2.72%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '83.40'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.81'
$ws.Range("E40").Value = '  +4.06%  '
$ws.Range("E41").Value = '  +0.20%  '
$ws.Range("E42").Value = '  +1.85%  '
$ws.Range("E43").Value = '  -0.56%  '
$ws.Range("E44").Value = '  -0.67%  '
$ws.Range("E45").Value = '  -1.00%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '5.95'
$ws.Range("E46").Value = '  +4.51%  '
$ws.Range("D47").Value = '1.946.39'
$ws.Range("E47").Value = '  +1.31%  '
$ws.Range("B48").Value = 'InjectiveProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '11.99'
$ws.Range("E48").Value = '  +3.39%  '
$ws.Range("B49").Value = 'Quant'
$ws.Range("C49").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '103.88'
$ws.Range("E49").Value = '  +6.70%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '49.78'
$ws.Range("E51").Value = '  -2.13%  '
